$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to reflect extended date range (01Mar-15May -> 15Feb-15May)
$ws.Name = "15Feb-15May"

# --- Column C ("Deaths") holds numeric-looking values stored as TEXT (shared strings). ---
# Force text number-format first so assigning a numeric-looking string keeps it text-typed,
# then restore the default style footprint once the text is stored.
$deathsRange = $ws.Range("C2:C19")
$deathsRange.NumberFormat = "@"

$ws.Range("C2").Value = "208320"
$ws.Range("C3").Value = "11719"
$ws.Range("C4").Value = "21703"
$ws.Range("C5").Value = "41198"
$ws.Range("C6").Value = "82259"
$ws.Range("C7").Value = "54392"
$ws.Range("C8").Value = "101568"
$ws.Range("C9").Value = "7478"
$ws.Range("C10").Value = "14234"
$ws.Range("C11").Value = "25314"
$ws.Range("C12").Value = "40136"
$ws.Range("C13").Value = "16522"
$ws.Range("C14").Value = "106754"
$ws.Range("C15").Value = "4242"
$ws.Range("C16").Value = "7469"
$ws.Range("C17").Value = "15886"
$ws.Range("C18").Value = "42128"
$ws.Range("C19").Value = "37870"

$deathsRange.Style = "Normal"

# --- Columns D:I hold true numeric Excess/ExcessPer stats; update each cell. ---
$ws.Range("D2").Value = 47490
$ws.Range("E2").Value = 43984
$ws.Range("F2").Value = 50362
$ws.Range("G2").Value = 29.5
$ws.Range("H2").Value = 26.8
$ws.Range("I2").Value = 31.9
$ws.Range("D3").Value = 1088
$ws.Range("E3").Value = -426
$ws.Range("F3").Value = 2209
$ws.Range("G3").Value = 10.2
$ws.Range("H3").Value = -3.5
$ws.Range("I3").Value = 23.2
$ws.Range("D4").Value = 4386
$ws.Range("E4").Value = 3286
$ws.Range("F4").Value = 5157
$ws.Range("G4").Value = 25.3
$ws.Range("H4").Value = 17.8
$ws.Range("I4").Value = 31.2
$ws.Range("D5").Value = 10614
$ws.Range("E5").Value = 8256
$ws.Range("F5").Value = 12524
$ws.Range("G5").Value = 34.7
$ws.Range("H5").Value = 25.1
$ws.Range("I5").Value = 43.7
$ws.Range("D6").Value = 20674
$ws.Range("E6").Value = 19377
$ws.Range("F6").Value = 21387
$ws.Range("G6").Value = 33.6
$ws.Range("H6").Value = 30.8
$ws.Range("I6").Value = 35.1
$ws.Range("D7").Value = 12114
$ws.Range("E7").Value = 10877
$ws.Range("F7").Value = 12858
$ws.Range("G7").Value = 28.7
$ws.Range("H7").Value = 25
$ws.Range("I7").Value = 31
$ws.Range("D8").Value = 24655
$ws.Range("E8").Value = 22604
$ws.Range("F8").Value = 26215
$ws.Range("G8").Value = 32.1
$ws.Range("H8").Value = 28.6
$ws.Range("I8").Value = 34.8
$ws.Range("D9").Value = 935
$ws.Range("E9").Value = 246
$ws.Range("F9").Value = 1342
$ws.Range("G9").Value = 14.3
$ws.Range("H9").Value = 3.4
$ws.Range("I9").Value = 21.9
$ws.Range("D10").Value = 3505
$ws.Range("E10").Value = 2646
$ws.Range("F10").Value = 4006
$ws.Range("G10").Value = 32.7
$ws.Range("H10").Value = 22.8
$ws.Range("I10").Value = 39.2
$ws.Range("D11").Value = 7219
$ws.Range("E11").Value = 5646
$ws.Range("F11").Value = 8364
$ws.Range("G11").Value = 39.9
$ws.Range("H11").Value = 28.7
$ws.Range("I11").Value = 49.3
$ws.Range("D12").Value = 10855
$ws.Range("E12").Value = 9410
$ws.Range("F12").Value = 11716
$ws.Range("G12").Value = 37.1
$ws.Range("H12").Value = 30.6
$ws.Range("I12").Value = 41.2
$ws.Range("D13").Value = 3623
$ws.Range("E13").Value = 2849
$ws.Range("F13").Value = 4055
$ws.Range("G13").Value = 28.1
$ws.Range("H13").Value = 20.8
$ws.Range("I13").Value = 32.5
$ws.Range("D14").Value = 23125
$ws.Range("E14").Value = 20997
$ws.Range("F14").Value = 24609
$ws.Range("G14").Value = 27.7
$ws.Range("H14").Value = 24.5
$ws.Range("I14").Value = 30
$ws.Range("D15").Value = 430
$ws.Range("E15").Value = -242
$ws.Range("F15").Value = 882
$ws.Range("G15").Value = 11.3
$ws.Range("H15").Value = -5.4
$ws.Range("I15").Value = 26.2
$ws.Range("D16").Value = 1154
$ws.Range("E16").Value = 739
$ws.Range("F16").Value = 1388
$ws.Range("G16").Value = 18.3
$ws.Range("H16").Value = 11
$ws.Range("I16").Value = 22.8
$ws.Range("D17").Value = 3676
$ws.Range("E17").Value = 2399
$ws.Range("F17").Value = 4438
$ws.Range("G17").Value = 30.1
$ws.Range("H17").Value = 17.8
$ws.Range("I17").Value = 38.8
$ws.Range("D18").Value = 10090
$ws.Range("E18").Value = 9181
$ws.Range("F18").Value = 10629
$ws.Range("G18").Value = 31.5
$ws.Range("H18").Value = 27.9
$ws.Range("I18").Value = 33.7
$ws.Range("D19").Value = 8806
$ws.Range("E19").Value = 7590
$ws.Range("F19").Value = 9536
$ws.Range("G19").Value = 30.3
$ws.Range("H19").Value = 25.1
$ws.Range("I19").Value = 33.7
